$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0.001
$ws.Range("K11").Value = 477
$ws.Range("L11").Value = 0.002385
